$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (values that Excel will not mis-detect as numbers)
$ws.Range("D2").Value = '27.099.54'
$ws.Range("E2").Value = '  -2.42%  '
$ws.Range("D3").Value = '1.560.51'
$ws.Range("E3").Value = '  -2.30%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  -1.41%  '
$ws.Range("E6").Value = '  -3.66%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -1.35%  '
$ws.Range("E9").Value = '  -3.07%  '
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").Value = '1.782.77'
$ws.Range("E12").Value = '  -2.25%  '
$ws.Range("D13").Value = '1.552.78'
$ws.Range("E13").Value = '  -2.74%  '
$ws.Range("E14").Value = '  -2.83%  '
$ws.Range("E15").Value = '  -3.18%  '
$ws.Range("E16").Value = '  -1.11%  '
$ws.Range("D17").Value = '27.095.31'
$ws.Range("E17").Value = '  -2.36%  '
$ws.Range("E18").Value = '  -2.53%  '
$ws.Range("E19").Value = '  -1.90%  '
$ws.Range("E20").Value = '  -2.26%  '
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("E22").Value = '  -1.48%  '
$ws.Range("E23").Value = '  -4.95%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("E25").Value = '  -1.73%  '
$ws.Range("E26").Value = '  -7.86%  '
$ws.Range("E27").Value = '  -1.97%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("E29").Value = '  -1.67%  '
$ws.Range("E30").Value = '  -1.88%  '
$ws.Range("E31").Value = '  -2.61%  '
$ws.Range("E32").Value = '  -2.45%  '
$ws.Range("D33").Value = '1.386.99'
$ws.Range("E33").Value = '  +0.66%  '
$ws.Range("E34").Value = '  -2.24%  '
$ws.Range("E35").Value = '  -0.44%  '
$ws.Range("E36").Value = '  -1.91%  '
$ws.Range("E37").Value = '  -4.09%  '
$ws.Range("E38").Value = '  -2.35%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("E39").Value = '  -3.83%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("E40").Value = '  -2.48%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  +2.09%  '
$ws.Range("E43").Value = '  +1.49%  '
$ws.Range("E44").Value = '  -2.06%  '
$ws.Range("E45").Value = '  -0.79%  '
$ws.Range("E46").Value = '  -0.19%  '
$ws.Range("D47").Value = '1.695.67'
$ws.Range("E47").Value = '  -2.19%  '
$ws.Range("E48").Value = '  -1.68%  '
$ws.Range("D49").Value = '0.0₇0988'
$ws.Range("E49").Value = '  -1.91%  '
$ws.Range("E50").Value = '  -1.02%  '
$ws.Range("E51").Value = '  -0.09%  '

# Numeric-looking text updates: force text format so Excel keeps the literal string
# (otherwise values like 0.810 or 206.25 would be auto-converted to numbers,
# losing trailing zeros / introducing float artifacts), then restore default style.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.485'
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0589'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0864'
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.516'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.33'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.939'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0165'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.517'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.810'
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.990'
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0492'
$ws.Range("D50").Style = "Normal"
